$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "41.199.38"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "2.175.76"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.27"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.41"
$ws.Range("E7").Value = "  -5.13%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  -5.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("E10").Value = "  -9.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  -3.55%  "
$ws.Range("E12").Value = "  -5.17%  "
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").Value = "2.498.71"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.95"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.810"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("D17").Value = "2.178.43"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "41.042.73"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").Value = "  -7.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.47"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("E22").Value = "  -10.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.29"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -7.16%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.93"
$ws.Range("E26").Value = "  -6.00%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.60"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.02"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.15"
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0773"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("E34").Value = "  -9.32%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -8.72%  "
$ws.Range("E37").Value = "  -4.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0287"
$ws.Range("E38").Value = "  -5.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.15"
$ws.Range("E39").Value = "  -6.86%  "
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.45"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "60.58"
$ws.Range("E42").Value = "  -7.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.192"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.35"
$ws.Range("E44").Value = "  -5.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0976"
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.40"
$ws.Range("E46").Value = "  -6.02%  "
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  -8.74%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").Value = "2.377.85"
$ws.Range("E51").Value = "  -2.06%  "
